$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2088.353
$ws.Range("I40").Value = 2099.9
$ws.Range("J40").Value = 2071.8572
$ws.Range("K40").Value = 2099.9
$ws.Range("L40").Value = 2071.8572
$ws.Range("M40").Value = -1924.9
$ws.Range("N40").Value = -2421.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3209.0986
$ws.Range("I32").Value = 3172.0598
$ws.Range("J32").Value = 3829.5
$ws.Range("K32").Value = 3172.0598
$ws.Range("L32").Value = 3829.5
$ws.Range("M32").Value = -2885.0598
$ws.Range("N32").Value = -4403.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3300.7805
$ws.Range("I61").Value = 3470.389
$ws.Range("J61").Value = 2079.6
$ws.Range("K61").Value = 3470.389
$ws.Range("L61").Value = 2079.6
$ws.Range("M61").Value = -3258.389
$ws.Range("N61").Value = -2503.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3300.7805
$ws.Range("I136").Value = 3470.389
$ws.Range("J136").Value = 2079.6
$ws.Range("K136").Value = 10411.167
$ws.Range("L136").Value = 6238.799999999999
$ws.Range("M136").Value = -7861.167000000001
$ws.Range("N136").Value = -11338.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 11296.333
$ws.Range("I134").Value = 11700.521
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 35101.563
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -32566.563
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 914.4666999999999
$ws.Range("I16").Value = 898.875
$ws.Range("J16").Value = 932.2857
$ws.Range("K16").Value = 898.875
$ws.Range("L16").Value = 932.2857
$ws.Range("M16").Value = -611.875
$ws.Range("N16").Value = -1506.2857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11115.5
$ws.Range("I31").Value = 5088.56
$ws.Range("J31").Value = 16696
$ws.Range("K31").Value = 5088.56
$ws.Range("L31").Value = 16696
$ws.Range("M31").Value = -4793.56
$ws.Range("N31").Value = -17286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 11115.5
$ws.Range("I34").Value = 5088.56
$ws.Range("J34").Value = 16696
$ws.Range("K34").Value = 5088.56
$ws.Range("L34").Value = 16696
$ws.Range("M34").Value = -4886.56
$ws.Range("N34").Value = -17100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 4500
$ws.Range("I41").Value = 4500
$ws.Range("J41").Value = 4500
$ws.Range("K41").Value = 4500
$ws.Range("L41").Value = 4500
$ws.Range("M41").Value = -4072
$ws.Range("N41").Value = -5356

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 8895
$ws.Range("I50").Value = 5900
$ws.Range("J50").Value = 9394.166999999999
$ws.Range("K50").Value = 5900
$ws.Range("L50").Value = 9394.166999999999
$ws.Range("M50").Value = -5275
$ws.Range("N50").Value = -10644.167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 10181.857
$ws.Range("I51").Value = 15000
$ws.Range("J51").Value = 9378.833000000001
$ws.Range("K51").Value = 15000
$ws.Range("L51").Value = 9378.833000000001
$ws.Range("M51").Value = -14264
$ws.Range("N51").Value = -10850.833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 15935.6
$ws.Range("J59").Value = 15935.6
$ws.Range("L59").Value = 15935.6
$ws.Range("N59").Value = -18225.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 8519
$ws.Range("J60").Value = 8519
$ws.Range("L60").Value = 8519
$ws.Range("N60").Value = -9541

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 10181.857
$ws.Range("I61").Value = 15000
$ws.Range("J61").Value = 9378.833000000001
$ws.Range("K61").Value = 15000
$ws.Range("L61").Value = 9378.833000000001
$ws.Range("M61").Value = -14652
$ws.Range("N61").Value = -10074.833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 17832.666
$ws.Range("J68").Value = 17832.666
$ws.Range("L68").Value = 17832.666
$ws.Range("N68").Value = -19330.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 17832.666
$ws.Range("J71").Value = 17832.666
$ws.Range("L71").Value = 53497.99800000001
$ws.Range("N71").Value = -60985.99800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 16571.5
$ws.Range("J74").Value = 16571.5
$ws.Range("L74").Value = 16571.5
$ws.Range("N74").Value = -18319.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 16571.5
$ws.Range("J77").Value = 16571.5
$ws.Range("L77").Value = 49714.5
$ws.Range("N77").Value = -58450.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 914.4666999999999
$ws.Range("I113").Value = 898.875
$ws.Range("J113").Value = 932.2857
$ws.Range("K113").Value = 898.875
$ws.Range("L113").Value = 932.2857
$ws.Range("M113").Value = 1271.125
$ws.Range("N113").Value = -5272.2857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4099243.8
$ws.Range("I134").Value = 4386708
$ws.Range("J134").Value = 2875
$ws.Range("K134").Value = 13160124
$ws.Range("L134").Value = 8625
$ws.Range("M134").Value = -13157589
$ws.Range("N134").Value = -13695

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 35712.47
$ws.Range("I70").Value = 51966.617
$ws.Range("J70").Value = 4681.8184
$ws.Range("K70").Value = 51966.617
$ws.Range("L70").Value = 4681.8184
$ws.Range("M70").Value = -51696.617
$ws.Range("N70").Value = -5221.8184

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 35712.47
$ws.Range("I73").Value = 51966.617
$ws.Range("J73").Value = 4681.8184
$ws.Range("K73").Value = 51966.617
$ws.Range("L73").Value = 4681.8184
$ws.Range("M73").Value = -51030.617
$ws.Range("N73").Value = -6553.8184

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1050
$ws.Range("I40").Value = 1050
$ws.Range("K40").Value = 1050
$ws.Range("M40").Value = -914

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H51").Value = 50042
$ws.Range("J51").Value = 50042
$ws.Range("L51").Value = 50042
$ws.Range("N51").Value = -50998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3567.5557
$ws.Range("I61").Value = 2833.3333
$ws.Range("J61").Value = 3934.6667
$ws.Range("K61").Value = 2833.3333
$ws.Range("L61").Value = 3934.6667
$ws.Range("M61").Value = -2631.3333
$ws.Range("N61").Value = -4338.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3567.5557
$ws.Range("I113").Value = 2833.3333
$ws.Range("J113").Value = 3934.6667
$ws.Range("K113").Value = 2833.3333
$ws.Range("L113").Value = 3934.6667
$ws.Range("M113").Value = -663.3332999999998
$ws.Range("N113").Value = -8274.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3881.82
$ws.Range("I136").Value = 4551.243
$ws.Range("J136").Value = 1976.5385
$ws.Range("K136").Value = 13653.729
$ws.Range("L136").Value = 5929.6155
$ws.Range("M136").Value = -11103.729
$ws.Range("N136").Value = -11029.6155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 14290.939
$ws.Range("I136").Value = 16723.812
$ws.Range("J136").Value = 2300.3572
$ws.Range("K136").Value = 50171.436
$ws.Range("L136").Value = 6901.071599999999
$ws.Range("M136").Value = -47621.436
$ws.Range("N136").Value = -12001.0716
